$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2312.5
$ws.Range("I2").Value = 2312.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2312.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2199.5
$ws.Range("N2").ClearContents()

$ws.Range("H6").Value = 433.96155
$ws.Range("I6").Value = 419.96
$ws.Range("K6").Value = 1259.88
$ws.Range("M6").Value = -1147.88

$ws.Range("H9").Value = 205.81818
$ws.Range("I9").Value = 234.66667
$ws.Range("K9").Value = 234.66667
$ws.Range("M9").Value = -65.66667000000001

$ws.Range("H19").Value = 1305.2778
$ws.Range("I19").Value = 682.1111
$ws.Range("J19").Value = 1928.4445
$ws.Range("K19").Value = 682.1111
$ws.Range("L19").Value = 1928.4445
$ws.Range("M19").Value = -507.1111
$ws.Range("N19").Value = -2278.4445

$ws.Range("H29").Value = 11284.143
$ws.Range("J29").Value = 11664.333
$ws.Range("L29").Value = 34992.999
$ws.Range("N29").Value = -35554.999

$ws.Range("H38").Value = 675.25
$ws.Range("I38").Value = 200.42857
$ws.Range("J38").Value = 3999
$ws.Range("K38").Value = 601.28571
$ws.Range("L38").Value = 11997
$ws.Range("M38").Value = -229.28571
$ws.Range("N38").Value = -12741

$ws.Range("H41").Value = 3321.5625
$ws.Range("I41").Value = 3387.6
$ws.Range("J41").Value = 3211.5
$ws.Range("K41").Value = 3387.6
$ws.Range("L41").Value = 3211.5
$ws.Range("M41").Value = -2947.6
$ws.Range("N41").Value = -4091.5

$ws.Range("H53").Value = 524
$ws.Range("I53").Value = 421.22223
$ws.Range("J53").Value = 656.1429000000001
$ws.Range("K53").Value = 421.22223
$ws.Range("L53").Value = 656.1429000000001
$ws.Range("M53").Value = 215.77777
$ws.Range("N53").Value = -1930.1429

$ws.Range("H86").Value = 695438.25
$ws.Range("I86").Value = 1677059.8
$ws.Range("J86").Value = 2528.9412
$ws.Range("K86").Value = 1677059.8
$ws.Range("L86").Value = 2528.9412
$ws.Range("M86").Value = -1675936.8
$ws.Range("N86").Value = -4774.9412

$ws.Range("H89").Value = 695438.25
$ws.Range("I89").Value = 1677059.8
$ws.Range("J89").Value = 2528.9412
$ws.Range("K89").Value = 8385299
$ws.Range("L89").Value = 12644.706
$ws.Range("M89").Value = -8379683
$ws.Range("N89").Value = -23876.706

$ws.Range("H112").Value = 1149.6052
$ws.Range("J112").Value = 1280.1613
$ws.Range("L112").Value = 3840.4839
$ws.Range("N112").Value = -6056.4839

$ws.Range("H113").Value = 62505800
$ws.Range("I113").Value = 25006224
$ws.Range("J113").Value = 100005380
$ws.Range("K113").Value = 25006224
$ws.Range("L113").Value = 100005380
$ws.Range("M113").Value = -25002970
$ws.Range("N113").Value = -100011888

$ws.Range("H125").Value = 3496
$ws.Range("J125").Value = 3496
$ws.Range("L125").Value = 31464
$ws.Range("N125").Value = -36384

$ws.Range("H132").Value = 2778.7896
$ws.Range("I132").Value = 2608.3333
$ws.Range("K132").Value = 7824.999899999999
$ws.Range("M132").Value = -5294.999899999999

$ws.Range("H137").Value = 3374.4468
$ws.Range("I137").Value = 2358.0322
$ws.Range("K137").Value = 7074.096600000001
$ws.Range("M137").Value = -4524.096600000001

$ws.Range("H138").Value = 3117.8042
$ws.Range("J138").Value = 3377.2168
$ws.Range("L138").Value = 10131.6504
$ws.Range("N138").Value = -20411.6504

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 3393.0557
$ws.Range("I141").Value = 3568.3635
$ws.Range("K141").Value = 10705.0905
$ws.Range("M141").Value = -5525.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 728168
$ws.Range("I16").Value = 872601.8
$ws.Range("J16").Value = 5999
$ws.Range("K16").Value = 872601.8
$ws.Range("L16").Value = 5999
$ws.Range("M16").Value = -872314.8
$ws.Range("N16").Value = -6573

$ws.Range("H32").Value = 11630461
$ws.Range("I32").Value = 12501545
$ws.Range("K32").Value = 12501545
$ws.Range("M32").Value = -12501258

$ws.Range("I61").Value = 100013170
$ws.Range("J61").Value = 100056.664
$ws.Range("K61").Value = 100013170
$ws.Range("L61").Value = 100056.664
$ws.Range("M61").Value = -100012958
$ws.Range("N61").Value = -100480.664

$ws.Range("H74").Value = 5688443
$ws.Range("I74").Value = 7577949
$ws.Range("J74").Value = 19926.455
$ws.Range("K74").Value = 7577949
$ws.Range("L74").Value = 19926.455
$ws.Range("M74").Value = -7577075
$ws.Range("N74").Value = -21674.455

$ws.Range("H77").Value = 5688443
$ws.Range("I77").Value = 7577949
$ws.Range("J77").Value = 19926.455
$ws.Range("K77").Value = 37889745
$ws.Range("L77").Value = 99632.27500000001
$ws.Range("M77").Value = -37885377
$ws.Range("N77").Value = -108368.275

$ws.Range("H94").Value = 47232
$ws.Range("J94").Value = 47232
$ws.Range("L94").Value = 47232
$ws.Range("N94").Value = -49034

$ws.Range("H97").Value = 1712.3125
$ws.Range("I97").Value = 1889.4286
$ws.Range("J97").Value = 472.5
$ws.Range("K97").Value = 1889.4286
$ws.Range("L97").Value = 472.5
$ws.Range("M97").Value = -1393.4286
$ws.Range("N97").Value = -1464.5

$ws.Range("H124").Value = 19760
$ws.Range("J124").Value = 19760
$ws.Range("L124").Value = 19760
$ws.Range("N124").Value = -29580

$ws.Range("H125").Value = 67568.336
$ws.Range("J125").Value = 67568.336
$ws.Range("L125").Value = 67568.336
$ws.Range("N125").Value = -77408.336

$ws.Range("H132").Value = 5117.2104
$ws.Range("I132").Value = 1890.3846
$ws.Range("K132").Value = 5671.1538
$ws.Range("M132").Value = -3141.1538

$ws.Range("H135").Value = 112500
$ws.Range("J135").Value = 112500
$ws.Range("L135").Value = 112500
$ws.Range("N135").Value = -122640

$ws.Range("I136").Value = 100013170
$ws.Range("J136").Value = 100056.664
$ws.Range("K136").Value = 300039510
$ws.Range("L136").Value = 300169.992
$ws.Range("M136").Value = -300036960
$ws.Range("N136").Value = -305269.992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2518.4375
$ws.Range("I20").Value = 2634.1724
$ws.Range("K20").Value = 2634.1724
$ws.Range("M20").Value = -2387.1724

$ws.Range("H86").Value = 2676.8948
$ws.Range("I86").Value = 1953.5555
$ws.Range("J86").Value = 3327.9
$ws.Range("K86").Value = 1953.5555
$ws.Range("L86").Value = 3327.9
$ws.Range("M86").Value = -830.5554999999999
$ws.Range("N86").Value = -5573.9

$ws.Range("H89").Value = 2676.8948
$ws.Range("I89").Value = 1953.5555
$ws.Range("J89").Value = 3327.9
$ws.Range("K89").Value = 9767.7775
$ws.Range("L89").Value = 16639.5
$ws.Range("M89").Value = -4151.7775
$ws.Range("N89").Value = -27871.5

$ws.Range("H94").Value = 747.4167
$ws.Range("I94").Value = 568
$ws.Range("K94").Value = 568
$ws.Range("M94").Value = -117

$ws.Range("H130").Value = 66662.336
$ws.Range("J130").Value = 66662.336
$ws.Range("L130").Value = 66662.336
$ws.Range("N130").Value = -76702.336

$ws.Range("H135").Value = 60606.062
$ws.Range("J135").Value = 60606.062
$ws.Range("L135").Value = 60606.062
$ws.Range("N135").Value = -70746.06200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 563832.8
$ws.Range("I31").Value = 9140.24
$ws.Range("K31").Value = 9140.24
$ws.Range("M31").Value = -8845.24

$ws.Range("H33").Value = 635.4286
$ws.Range("I33").Value = 635.4286
$ws.Range("K33").Value = 635.4286
$ws.Range("M33").Value = -256.4286

$ws.Range("H34").Value = 563832.8
$ws.Range("I34").Value = 9140.24
$ws.Range("K34").Value = 9140.24
$ws.Range("M34").Value = -8938.24

$ws.Range("H58").Value = 5817.643
$ws.Range("I58").Value = 1589.4445
$ws.Range("J58").Value = 13428.4
$ws.Range("K58").Value = 1589.4445
$ws.Range("L58").Value = 13428.4
$ws.Range("M58").Value = -1386.4445
$ws.Range("N58").Value = -13834.4

$ws.Range("H105").Value = 2458.6155
$ws.Range("I105").Value = 2501
$ws.Range("J105").Value = 2280.6
$ws.Range("K105").Value = 2501
$ws.Range("L105").Value = 2280.6
$ws.Range("M105").Value = -754
$ws.Range("N105").Value = -5774.6

$ws.Range("H107").Value = 1638.8
$ws.Range("I107").Value = 1703
$ws.Range("J107").Value = 1604.2307
$ws.Range("K107").Value = 1703
$ws.Range("L107").Value = 1604.2307
$ws.Range("M107").Value = 217
$ws.Range("N107").Value = -5444.2307

$ws.Range("H136").Value = 5817.643
$ws.Range("I136").Value = 1589.4445
$ws.Range("J136").Value = 13428.4
$ws.Range("K136").Value = 4768.333500000001
$ws.Range("L136").Value = 40285.2
$ws.Range("M136").Value = -2218.333500000001
$ws.Range("N136").Value = -45385.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1857.8334
$ws.Range("I13").Value = 81.666664
$ws.Range("J13").Value = 3634
$ws.Range("K13").Value = 244.999992
$ws.Range("L13").Value = 10902
$ws.Range("M13").Value = -76.99999199999999
$ws.Range("N13").Value = -11238

$ws.Range("H44").Value = 17121.9
$ws.Range("I44").Value = 21388
$ws.Range("K44").Value = 64164
$ws.Range("M44").Value = -63766

$ws.Range("H46").Value = 874.25
$ws.Range("I46").Value = 499.33334
$ws.Range("K46").Value = 1498.00002
$ws.Range("M46").Value = -1407.00002

$ws.Range("H92").Value = 10000002
$ws.Range("J92").Value = 10000003
$ws.Range("L92").Value = 30000009
$ws.Range("N92").Value = -30002505

$ws.Range("H122").Value = 2216.25
$ws.Range("I122").Value = 972
$ws.Range("K122").Value = 8748
$ws.Range("M122").Value = -6298

$ws.Range("H126").Value = 6139.25
$ws.Range("I126").Value = 4566.364
$ws.Range("K126").Value = 13699.092
$ws.Range("M126").Value = -8759.091999999999

$ws.Range("H132").Value = 2240.6428
$ws.Range("I132").Value = 2000.1666
$ws.Range("J132").Value = 2421
$ws.Range("K132").Value = 18001.4994
$ws.Range("L132").Value = 21789
$ws.Range("M132").Value = -15471.4994
$ws.Range("N132").Value = -26849

$ws.Range("H139").Value = 2736.45
$ws.Range("I139").Value = 2528.625
$ws.Range("K139").Value = 7585.875
$ws.Range("M139").Value = -2445.875

$ws.Range("H141").Value = 191351.38
$ws.Range("I141").Value = 191351.38
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 574054.14
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -568874.14
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 704
$ws.Range("I23").Value = 212
$ws.Range("J23").Value = 950
$ws.Range("K23").Value = 212
$ws.Range("L23").Value = 950
$ws.Range("M23").Value = 11
$ws.Range("N23").Value = -1396

$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -4710
$ws.Range("N29").ClearContents()

$ws.Range("H70").Value = 5271.7144
$ws.Range("I70").Value = 5480.4
$ws.Range("J70").Value = 4750
$ws.Range("K70").Value = 5480.4
$ws.Range("L70").Value = 4750
$ws.Range("M70").Value = -5210.4
$ws.Range("N70").Value = -5290

$ws.Range("H73").Value = 5271.7144
$ws.Range("I73").Value = 5480.4
$ws.Range("J73").Value = 4750
$ws.Range("K73").Value = 5480.4
$ws.Range("L73").Value = 4750
$ws.Range("M73").Value = -4544.4
$ws.Range("N73").Value = -6622

$ws.Range("H80").Value = 600
$ws.Range("I80").Value = 400
$ws.Range("J80").Value = 640
$ws.Range("K80").Value = 400
$ws.Range("L80").Value = 640
$ws.Range("M80").Value = 598
$ws.Range("N80").Value = -2636

$ws.Range("H83").Value = 600
$ws.Range("I83").Value = 400
$ws.Range("J83").Value = 640
$ws.Range("K83").Value = 2000
$ws.Range("L83").Value = 3200
$ws.Range("M83").Value = 2992
$ws.Range("N83").Value = -13184

$ws.Range("H93").Value = 59999.5
$ws.Range("J93").Value = 59999.5
$ws.Range("L93").Value = 59999.5
$ws.Range("N93").Value = -63743.5

$ws.Range("H102").Value = 9562.727999999999
$ws.Range("I102").Value = 7798.8887
$ws.Range("K102").Value = 7798.8887
$ws.Range("M102").Value = -6176.8887

$ws.Range("H107").Value = 294.375
$ws.Range("I107").Value = 279.2857
$ws.Range("K107").Value = 279.2857
$ws.Range("M107").Value = 1640.7143

$ws.Range("H122").Value = 3621.7778
$ws.Range("I122").Value = 2798.1667
$ws.Range("J122").Value = 5269
$ws.Range("K122").Value = 8394.500100000001
$ws.Range("L122").Value = 15807
$ws.Range("M122").Value = -5944.500100000001
$ws.Range("N122").Value = -20707

$ws.Range("H132").Value = 50008224
$ws.Range("I132").Value = 90911336
$ws.Range("K132").Value = 272734008
$ws.Range("M132").Value = -272731478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1965.3334
$ws.Range("I16").Value = 1618.2667
$ws.Range("K16").Value = 1618.2667
$ws.Range("M16").Value = -1448.2667

$ws.Range("H23").Value = 3999.875
$ws.Range("I23").Value = 999.5
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 999.5
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -769.5
$ws.Range("N23").Value = -5460

$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 18000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 18000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -18580

$ws.Range("H82").Value = 1308
$ws.Range("I82").Value = 1197.625
$ws.Range("J82").Value = 1749.5
$ws.Range("K82").Value = 1197.625
$ws.Range("L82").Value = 1749.5
$ws.Range("M82").Value = -836.625
$ws.Range("N82").Value = -2471.5

$ws.Range("H85").Value = 1308
$ws.Range("I85").Value = 1197.625
$ws.Range("J85").Value = 1749.5
$ws.Range("K85").Value = 1197.625
$ws.Range("L85").Value = 1749.5
$ws.Range("M85").Value = 50.375
$ws.Range("N85").Value = -4245.5

$ws.Range("H127").Value = 155943
$ws.Range("J127").Value = 155943
$ws.Range("L127").Value = 155943
$ws.Range("N127").Value = -165863

$ws.Range("H132").Value = 558458.75
$ws.Range("J132").Value = 3336670
$ws.Range("L132").Value = 10010010
$ws.Range("N132").Value = -10015070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21138

$ws.Range("H50").Value = 10000
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = -9369
$ws.Range("N50").Value = -11262

$ws.Range("H70").Value = 65275.668
$ws.Range("J70").Value = 65275.668
$ws.Range("L70").Value = 65275.668
$ws.Range("N70").Value = -65905.66800000001

$ws.Range("H73").Value = 65275.668
$ws.Range("J73").Value = 65275.668
$ws.Range("L73").Value = 65275.668
$ws.Range("N73").Value = -67459.66800000001

$ws.Range("H74").Value = 27500
$ws.Range("J74").Value = 27500
$ws.Range("L74").Value = 27500
$ws.Range("N74").Value = -29372

$ws.Range("H77").Value = 27500
$ws.Range("J77").Value = 27500
$ws.Range("L77").Value = 82500
$ws.Range("N77").Value = -91860

$ws.Range("H96").Value = 5000.4
$ws.Range("I96").Value = 2001.5
$ws.Range("K96").Value = 2001.5
$ws.Range("M96").Value = -628.5

$ws.Range("H100").Value = 840.36365
$ws.Range("I100").Value = 966.875
$ws.Range("K100").Value = 1933.75
$ws.Range("M100").Value = -1392.75

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H132").Value = 1891.4286
$ws.Range("I132").Value = 1857.7778
$ws.Range("K132").Value = 5573.3334
$ws.Range("M132").Value = -3043.3334
